# feat: add 2022-Q4 data
#
# 1. Insert a new "2022-Q4" row into the "总计" (summary) sheet.
# 2. Insert a new "2022-Q4" worksheet (with its fund-holding detail data)
#    positioned right after "总计" and before "2022-Q3".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Part 1: "总计" summary sheet - insert new row 2 for 2022-Q4
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

# Push existing data rows down by inserting a fresh row at position 2.
$summary.Rows.Item(2).Insert()

# Re-number the "index" column (A) for the rows that shifted down.
for ($r = 3; $r -le 7; $r++) {
    $summary.Cells.Item($r, 1).Value = $r - 2
}

# Fill in the new 2022-Q4 summary row.
$summary.Cells.Item(2, 1).Value = 0
$summary.Cells.Item(2, 2).Value = "2022-Q4"
$summary.Cells.Item(2, 3).Value = 4
$summary.Cells.Item(2, 4).Value = 3.34

# The inserted row copied row 1's formatting - restore column A's "index"
# style (matches the other data rows) and strip the stray formatting that
# landed on B:D so they look like plain data cells again.
$summary.Cells.Item(3, 1).Copy()
$summary.Cells.Item(2, 1).PasteSpecial(-4122)
$summary.Range("B2:D2").ClearFormats()

# ---------------------------------------------------------------------
# Part 2: new "2022-Q4" worksheet with fund holdings detail
# ---------------------------------------------------------------------
# Clone the "2022-Q3" sheet so the new sheet starts with identical
# formatting/layout, dropped right before it in tab order.
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3)
$q4 = $wb.Worksheets.Item("2022-Q3 (2)")
$q4.Name = "2022-Q4"

# 2022-Q3 has 13 data rows (rows 2-14); 2022-Q4 only needs 4 (rows 2-5) -
# drop the extra rows inherited from the copy.
$q4.Rows("6:14").Delete()

function Set-TextCell($cell, $text) {
    # Force text storage (so things like fund codes keep leading zeros)
    # without leaving a lingering quote-prefix style on the cell.
    $cell.Value = "'" + $text
    $cell.ClearFormats()
}

# Row 2: 012930
Set-TextCell $q4.Cells.Item(2, 2) "012930"
Set-TextCell $q4.Cells.Item(2, 3) "中庚价值先锋股票"
Set-TextCell $q4.Cells.Item(2, 4) "68.71"
Set-TextCell $q4.Cells.Item(2, 5) "94.78"
Set-TextCell $q4.Cells.Item(2, 6) "4.21"
Set-TextCell $q4.Cells.Item(2, 7) "2.8927"
$q4.Cells.Item(2, 8).Value = 10

# Row 3: 001468
Set-TextCell $q4.Cells.Item(3, 2) "001468"
Set-TextCell $q4.Cells.Item(3, 3) "广发改革先锋灵活配置混合"
Set-TextCell $q4.Cells.Item(3, 4) "5.72"
Set-TextCell $q4.Cells.Item(3, 5) "93.29"
Set-TextCell $q4.Cells.Item(3, 6) "4.86"
Set-TextCell $q4.Cells.Item(3, 7) "0.2780"
$q4.Cells.Item(3, 8).Value = 2

# Row 4: 014062
Set-TextCell $q4.Cells.Item(4, 2) "014062"
Set-TextCell $q4.Cells.Item(4, 3) "景顺长城专精特新量化优选股票A"
Set-TextCell $q4.Cells.Item(4, 4) "7.66"
Set-TextCell $q4.Cells.Item(4, 5) "91.15"
Set-TextCell $q4.Cells.Item(4, 6) "1.44"
Set-TextCell $q4.Cells.Item(4, 7) "0.1103"
$q4.Cells.Item(4, 8).Value = 10

# Row 5: 014063
Set-TextCell $q4.Cells.Item(5, 2) "014063"
Set-TextCell $q4.Cells.Item(5, 3) "景顺长城专精特新量化优选股票C"
Set-TextCell $q4.Cells.Item(5, 4) "3.88"
Set-TextCell $q4.Cells.Item(5, 5) "91.15"
Set-TextCell $q4.Cells.Item(5, 6) "1.44"
Set-TextCell $q4.Cells.Item(5, 7) "0.0559"
$q4.Cells.Item(5, 8).Value = 10

$wb.Worksheets.Item("总计").Activate()
$wb.Worksheets.Item("总计").Range("A1").Select()
